$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5052.7646
$ws.Range("I74").Value = 4845.5386
$ws.Range("K74").Value = 4845.5386
$ws.Range("M74").Value = -3909.5386

$ws.Range("H77").Value = 5052.7646
$ws.Range("I77").Value = 4845.5386
$ws.Range("K77").Value = 24227.693
$ws.Range("M77").Value = -19547.693

$ws.Range("H129").Value = 1668.25
$ws.Range("I129").Value = 1668.25
$ws.Range("K129").Value = 5004.75
$ws.Range("M129").Value = -4.75

$ws.Range("H131").Value = 25001776
$ws.Range("I131").Value = 33334000
$ws.Range("J131").Value = 5105
$ws.Range("K131").Value = 100002000
$ws.Range("L131").Value = 15315
$ws.Range("M131").Value = -99996960
$ws.Range("N131").Value = -25395

$ws.Range("H132").Value = 2981.3
$ws.Range("I132").Value = 2201.6667
$ws.Range("K132").Value = 6605.000100000001
$ws.Range("M132").Value = -4075.000100000001

$ws.Range("H135").Value = 1283.7333
$ws.Range("I135").Value = 1131.9166
$ws.Range("K135").Value = 10187.2494
$ws.Range("M135").Value = -7652.249400000001

$ws.Range("H138").Value = 2320.7368
$ws.Range("I138").Value = 2160.2
$ws.Range("J138").Value = 2499.111
$ws.Range("K138").Value = 6480.599999999999
$ws.Range("L138").Value = 7497.333
$ws.Range("M138").Value = -1340.599999999999
$ws.Range("N138").Value = -17777.333

$ws.Range("H140").Value = 80764.37
$ws.Range("J140").Value = 80764.37
$ws.Range("L140").Value = 80764.37
$ws.Range("N140").Value = -91124.37

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 624.25
$ws.Range("I2").Value = 582.3333
$ws.Range("K2").Value = 582.3333
$ws.Range("M2").Value = -469.3333

$ws.Range("H61").Value = 43352.125
$ws.Range("I61").Value = 1572.6
$ws.Range("K61").Value = 1572.6
$ws.Range("M61").Value = -1360.6

$ws.Range("H102").Value = 57869
$ws.Range("I102").Value = 69377.92999999999
$ws.Range("J102").Value = 23342.2
$ws.Range("K102").Value = 69377.92999999999
$ws.Range("L102").Value = 23342.2
$ws.Range("M102").Value = -67755.92999999999
$ws.Range("N102").Value = -26586.2

$ws.Range("H116").Value = 624.25
$ws.Range("I116").Value = 582.3333
$ws.Range("K116").Value = 582.3333
$ws.Range("M116").Value = 1711.6667

$ws.Range("H130").Value = 19696.666
$ws.Range("J130").Value = 19696.666
$ws.Range("L130").Value = 19696.666
$ws.Range("N130").Value = -29736.666

$ws.Range("H132").Value = 2124.0881
$ws.Range("I132").Value = 1984.3704
$ws.Range("K132").Value = 5953.1112
$ws.Range("M132").Value = -3423.1112

$ws.Range("H136").Value = 43352.125
$ws.Range("I136").Value = 1572.6
$ws.Range("K136").Value = 4717.799999999999
$ws.Range("M136").Value = -2167.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 624.25
$ws.Range("I3").Value = 582.3333
$ws.Range("K3").Value = 582.3333
$ws.Range("M3").Value = -468.3333

$ws.Range("H132").Value = 30734.285
$ws.Range("J132").Value = 30734.285
$ws.Range("L132").Value = 30734.285
$ws.Range("N132").Value = -40854.285

$ws.Range("H134").Value = 5863.25
$ws.Range("I134").Value = 3579.0625
$ws.Range("K134").Value = 10737.1875
$ws.Range("M134").Value = -8202.1875

$ws.Range("H135").Value = 86523.336
$ws.Range("J135").Value = 86523.336
$ws.Range("L135").Value = 86523.336
$ws.Range("N135").Value = -96663.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3407.125
$ws.Range("I62").Value = 3525
$ws.Range("J62").Value = 3289.25
$ws.Range("K62").Value = 3525
$ws.Range("L62").Value = 3289.25
$ws.Range("M62").Value = -2901
$ws.Range("N62").Value = -4537.25

$ws.Range("H65").Value = 3407.125
$ws.Range("I65").Value = 3525
$ws.Range("J65").Value = 3289.25
$ws.Range("K65").Value = 17625
$ws.Range("L65").Value = 16446.25
$ws.Range("M65").Value = -14505
$ws.Range("N65").Value = -22686.25

$ws.Range("H122").Value = 2822
$ws.Range("I122").Value = 2389.9167
$ws.Range("J122").Value = 3859
$ws.Range("K122").Value = 7169.750100000001
$ws.Range("L122").Value = 11577
$ws.Range("M122").Value = -4719.750100000001
$ws.Range("N122").Value = -16477

$ws.Range("H134").Value = 2589488
$ws.Range("I134").Value = 3404236.2
$ws.Range("K134").Value = 10212708.6
$ws.Range("M134").Value = -10210173.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2704.7334
$ws.Range("I34").Value = 33.4
$ws.Range("J34").Value = 4040.4
$ws.Range("K34").Value = 100.2
$ws.Range("L34").Value = 12121.2
$ws.Range("M34").Value = -16.19999999999999
$ws.Range("N34").Value = -12289.2

$ws.Range("H56").Value = 5778.6665
$ws.Range("I56").Value = 5778.6665
$ws.Range("K56").Value = 5778.6665
$ws.Range("M56").Value = -5248.6665

$ws.Range("H107").Value = 773.1429000000001
$ws.Range("J107").Value = 478
$ws.Range("L107").Value = 1434
$ws.Range("N107").Value = -5274

$ws.Range("H110").Value = 7400
$ws.Range("I110").Value = 7400
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 22200
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -18110
$ws.Range("N110").ClearContents()

$ws.Range("H113").Value = 2860811
$ws.Range("J113").Value = 3740649
$ws.Range("L113").Value = 11221947
$ws.Range("N113").Value = -11226287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3250
$ws.Range("I80").Value = 3250
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3250
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2252
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 3250
$ws.Range("I83").Value = 3250
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 16250
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -11258
$ws.Range("N83").ClearContents()

$ws.Range("H132").Value = 3746.4119
$ws.Range("I132").Value = 3057.9565
$ws.Range("J132").Value = 5185.909
$ws.Range("K132").Value = 9173.869499999999
$ws.Range("L132").Value = 15557.727
$ws.Range("M132").Value = -6643.869499999999
$ws.Range("N132").Value = -20617.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7073941
$ws.Range("I40").Value = 3829.875
$ws.Range("J40").Value = 25927570
$ws.Range("K40").Value = 3829.875
$ws.Range("L40").Value = 25927570
$ws.Range("M40").Value = -3693.875
$ws.Range("N40").Value = -25927842

$ws.Range("H46").Value = 2428
$ws.Range("I46").Value = 1273.8334
$ws.Range("J46").Value = 3197.4443
$ws.Range("K46").Value = 1273.8334
$ws.Range("L46").Value = 3197.4443
$ws.Range("M46").Value = -1085.8334
$ws.Range("N46").Value = -3573.4443

$ws.Range("H54").Value = 38826.8
$ws.Range("J54").Value = 38826.8
$ws.Range("L54").Value = 38826.8
$ws.Range("N54").Value = -40114.8

$ws.Range("H100").Value = 4610.067
$ws.Range("J100").Value = 5300
$ws.Range("L100").Value = 5300
$ws.Range("N100").Value = -6382

$ws.Range("H136").Value = 3823.158
$ws.Range("I136").Value = 5743.125
$ws.Range("K136").Value = 17229.375
$ws.Range("M136").Value = -14679.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws.Range("H113").Value = 1980.2222
$ws.Range("I113").Value = 1267.3334
$ws.Range("J113").Value = 2336.6667
$ws.Range("K113").Value = 3802.0002
$ws.Range("L113").Value = 7010.000100000001
$ws.Range("M113").Value = -1632.0002
$ws.Range("N113").Value = -11350.0001

$ws.Range("H136").Value = 1629.6
$ws.Range("I136").Value = 1629.6
$ws.Range("K136").Value = 4888.799999999999
$ws.Range("M136").Value = -2338.799999999999
